$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 06:46:50'
$ws.Range("A3").Value = 'Total filas: 36'
$ws.Cells.Item(16, 1).Value = '06:46:50'
$ws.Cells.Item(16, 2).Value = '06:56'
$ws.Cells.Item(16, 3).Value = '215A_EL PATO'
$ws.Cells.Item(16, 4).Value = 10
$ws.Cells.Item(16, 5).Value = 'LP1912'
$ws.Cells.Item(17, 1).Value = '05:57:13'
$ws.Cells.Item(17, 2).Value = '06:57'
$ws.Cells.Item(17, 3).Value = '215A_EL PATO'
$ws.Cells.Item(17, 4).Value = 60
$ws.Cells.Item(17, 5).Value = 'LP1912'
$ws.Cells.Item(18, 1).Value = '05:57:13'
$ws.Cells.Item(18, 2).Value = '06:59'
$ws.Cells.Item(18, 3).Value = '225_GOMEZ'
$ws.Cells.Item(18, 4).Value = 62
$ws.Cells.Item(18, 5).Value = 'LP1912'
$ws.Cells.Item(19, 1).Value = '06:17:28'
$ws.Cells.Item(19, 2).Value = '07:15'
$ws.Cells.Item(19, 3).Value = '215C_EL PATO'
$ws.Cells.Item(19, 4).Value = 58
$ws.Cells.Item(19, 5).Value = 'LP1912'
$ws.Cells.Item(20, 1).Value = '05:57:13'
$ws.Cells.Item(20, 2).Value = '07:16'
$ws.Cells.Item(20, 3).Value = '215C_EL PATO'
$ws.Cells.Item(20, 4).Value = 79
$ws.Cells.Item(20, 5).Value = 'LP1912'
$ws.Cells.Item(21, 1).Value = '05:57:13'
$ws.Cells.Item(21, 2).Value = '07:19'
$ws.Cells.Item(21, 3).Value = '14_ABASTO'
$ws.Cells.Item(21, 4).Value = 82
$ws.Cells.Item(21, 5).Value = 'LP1912'
$ws.Cells.Item(22, 1).Value = '06:46:50'
$ws.Cells.Item(22, 2).Value = '07:20'
$ws.Cells.Item(22, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(22, 4).Value = 34
$ws.Cells.Item(22, 5).Value = 'LP1912'
$ws.Cells.Item(23, 1).Value = '05:57:13'
$ws.Cells.Item(23, 2).Value = '07:21'
$ws.Cells.Item(23, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(23, 4).Value = 84
$ws.Cells.Item(23, 5).Value = 'LP1912'
$ws.Cells.Item(24, 1).Value = '06:17:28'
$ws.Cells.Item(24, 2).Value = '07:21'
$ws.Cells.Item(24, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(24, 4).Value = 64
$ws.Cells.Item(24, 5).Value = 'LP1912'
$ws.Cells.Item(25, 1).Value = '05:57:13'
$ws.Cells.Item(25, 2).Value = '07:29'
$ws.Cells.Item(25, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(25, 4).Value = 92
$ws.Cells.Item(25, 5).Value = 'LP1912'
$ws.Cells.Item(26, 1).Value = '05:57:13'
$ws.Cells.Item(26, 2).Value = '07:35'
$ws.Cells.Item(26, 3).Value = '10_OLMOS'
$ws.Cells.Item(26, 4).Value = 98
$ws.Cells.Item(26, 5).Value = 'LP1912'
$ws.Cells.Item(27, 1).Value = '06:17:28'
$ws.Cells.Item(27, 2).Value = '07:36'
$ws.Cells.Item(27, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(27, 4).Value = 79
$ws.Cells.Item(27, 5).Value = 'LP1912'
$ws.Cells.Item(28, 1).Value = '05:57:13'
$ws.Cells.Item(28, 2).Value = '07:37'
$ws.Cells.Item(28, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(28, 4).Value = 100
$ws.Cells.Item(28, 5).Value = 'LP1912'
$ws.Cells.Item(29, 1).Value = '06:46:50'
$ws.Cells.Item(29, 2).Value = '07:43'
$ws.Cells.Item(29, 3).Value = '215A_EL PATO'
$ws.Cells.Item(29, 4).Value = 57
$ws.Cells.Item(29, 5).Value = 'LP1912'
$ws.Cells.Item(30, 1).Value = '06:35:22'
$ws.Cells.Item(30, 2).Value = '07:44'
$ws.Cells.Item(30, 3).Value = '215A_EL PATO'
$ws.Cells.Item(30, 4).Value = 69
$ws.Cells.Item(30, 5).Value = 'LP1912'
$ws.Cells.Item(31, 1).Value = '05:57:13'
$ws.Cells.Item(31, 2).Value = '07:55'
$ws.Cells.Item(31, 3).Value = '14_ABASTO'
$ws.Cells.Item(31, 4).Value = 118
$ws.Cells.Item(31, 5).Value = 'LP1912'
$ws.Cells.Item(32, 1).Value = '06:46:50'
$ws.Cells.Item(32, 2).Value = '08:00'
$ws.Cells.Item(32, 3).Value = '17_ROMERO'
$ws.Cells.Item(32, 4).Value = 103
$ws.Cells.Item(32, 5).Value = 'LP1912'
$ws.Cells.Item(33, 1).Value = '06:46:50'
$ws.Cells.Item(33, 2).Value = '08:00'
$ws.Cells.Item(33, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(33, 4).Value = 74
$ws.Cells.Item(33, 5).Value = 'LP1912'
$ws.Cells.Item(34, 1).Value = '06:17:28'
$ws.Cells.Item(34, 2).Value = '08:01'
$ws.Cells.Item(34, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(34, 4).Value = 104
$ws.Cells.Item(34, 5).Value = 'LP1912'
$ws.Cells.Item(35, 1).Value = '06:46:50'
$ws.Cells.Item(35, 2).Value = '08:06'
$ws.Cells.Item(35, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(35, 4).Value = 91
$ws.Cells.Item(35, 5).Value = 'LP1912'
$ws.Cells.Item(36, 1).Value = '06:17:28'
$ws.Cells.Item(36, 2).Value = '08:11'
$ws.Cells.Item(36, 3).Value = '10_OLMOS'
$ws.Cells.Item(36, 4).Value = 114
$ws.Cells.Item(36, 5).Value = 'LP1912'
$ws.Cells.Item(37, 1).Value = '06:17:28'
$ws.Cells.Item(37, 2).Value = '08:13'
$ws.Cells.Item(37, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(37, 4).Value = 116
$ws.Cells.Item(37, 5).Value = 'LP1912'
$ws.Cells.Item(38, 1).Value = '06:35:22'
$ws.Cells.Item(38, 2).Value = '08:29'
$ws.Cells.Item(38, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(38, 4).Value = 114
$ws.Cells.Item(38, 5).Value = 'LP1912'
$ws.Cells.Item(39, 1).Value = '06:35:22'
$ws.Cells.Item(39, 2).Value = '08:29'
$ws.Cells.Item(39, 3).Value = '15_ABASTO'
$ws.Cells.Item(39, 4).Value = 114
$ws.Cells.Item(39, 5).Value = 'LP1912'
$ws.Cells.Item(40, 1).Value = '06:46:50'
$ws.Cells.Item(40, 2).Value = '08:41'
$ws.Cells.Item(40, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(40, 4).Value = 115
$ws.Cells.Item(40, 5).Value = 'LP1912'
$ws.Cells.Item(41, 1).Value = '06:46:50'
$ws.Cells.Item(41, 2).Value = '08:43'
$ws.Cells.Item(41, 3).Value = '215C_EL PATO'
$ws.Cells.Item(41, 4).Value = 117
$ws.Cells.Item(41, 5).Value = 'LP1912'

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 06:46:50'
$ws.Range("A3").Value = 'Total filas: 9'
$ws.Cells.Item(8, 1).Value = '06:46:50'
$ws.Cells.Item(8, 2).Value = '06:56'
$ws.Cells.Item(8, 3).Value = '215A_EL PATO'
$ws.Cells.Item(8, 4).Value = 10
$ws.Cells.Item(8, 5).Value = 'LP1912'
$ws.Cells.Item(9, 1).Value = '05:57:13'
$ws.Cells.Item(9, 2).Value = '06:57'
$ws.Cells.Item(9, 3).Value = '215A_EL PATO'
$ws.Cells.Item(9, 4).Value = 60
$ws.Cells.Item(9, 5).Value = 'LP1912'
$ws.Cells.Item(10, 1).Value = '06:17:28'
$ws.Cells.Item(10, 2).Value = '07:15'
$ws.Cells.Item(10, 3).Value = '215C_EL PATO'
$ws.Cells.Item(10, 4).Value = 58
$ws.Cells.Item(10, 5).Value = 'LP1912'
$ws.Cells.Item(11, 1).Value = '05:57:13'
$ws.Cells.Item(11, 2).Value = '07:16'
$ws.Cells.Item(11, 3).Value = '215C_EL PATO'
$ws.Cells.Item(11, 4).Value = 79
$ws.Cells.Item(11, 5).Value = 'LP1912'
$ws.Cells.Item(12, 1).Value = '06:46:50'
$ws.Cells.Item(12, 2).Value = '07:43'
$ws.Cells.Item(12, 3).Value = '215A_EL PATO'
$ws.Cells.Item(12, 4).Value = 57
$ws.Cells.Item(12, 5).Value = 'LP1912'
$ws.Cells.Item(13, 1).Value = '06:35:22'
$ws.Cells.Item(13, 2).Value = '07:44'
$ws.Cells.Item(13, 3).Value = '215A_EL PATO'
$ws.Cells.Item(13, 4).Value = 69
$ws.Cells.Item(13, 5).Value = 'LP1912'
$ws.Cells.Item(14, 1).Value = '06:46:50'
$ws.Cells.Item(14, 2).Value = '08:43'
$ws.Cells.Item(14, 3).Value = '215C_EL PATO'
$ws.Cells.Item(14, 4).Value = 117
$ws.Cells.Item(14, 5).Value = 'LP1912'

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 06:46:50'
$ws.Range("A3").Value = 'Total filas: 3'
$ws.Cells.Item(6, 1).Value = '06:46:50'
$ws.Cells.Item(6, 2).Value = '07:42'
$ws.Cells.Item(6, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(6, 4).Value = 56
$ws.Cells.Item(6, 5).Value = 'L6173'
$ws.Cells.Item(7, 1).Value = '05:57:13'
$ws.Cells.Item(7, 2).Value = '07:43'
$ws.Cells.Item(7, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(7, 4).Value = 106
$ws.Cells.Item(7, 5).Value = 'LP1912'
$ws.Cells.Item(8, 1).Value = '06:46:50'
$ws.Cells.Item(8, 2).Value = '08:35'
$ws.Cells.Item(8, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(8, 4).Value = 109
$ws.Cells.Item(8, 5).Value = 'L6173'
